# "Form submit data driven and parameters"
#
# Inserts three new worksheets (FormSubmit, FormSubmitByParameter,
# FormSubmitByDataProvider) right after "Select Input" and before the
# existing "Radio Buttons Demo" / "Simple Form Demo" / "RowColumnTable"
# sheets, populates them, and makes the last of the new sheets the active
# tab.

$wb = $excel.ActiveWorkbook
$selectInput = $wb.Worksheets.Item(1)

# --- Create the three new sheets, in order, right after "Select Input" ---
$formSubmit = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $selectInput)
$formSubmit.Name = "FormSubmit"

$formSubmitByParameter = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $formSubmit)
$formSubmitByParameter.Name = "FormSubmitByParameter"

$formSubmitByDataProvider = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $formSubmitByParameter)
$formSubmitByDataProvider.Name = "FormSubmitByDataProvider"

# --- FormSubmit: data row typed first (so the shared-string table picks up
#     the values before the headers), then the header row, then the zip ---
$formSubmit.Activate()
$formSubmit.Range("A2").Value = "Sneha"
$formSubmit.Range("C2").Value = "Sneha95"
$formSubmit.Range("B2").Value = "George"
$formSubmit.Range("D2").Value = "Kottayam"
$formSubmit.Range("E2").Value = "Kerala"

$formSubmit.Range("A1").Value = "First Name"
$formSubmit.Range("B1").Value = "Last Name"
$formSubmit.Range("C1").Value = "UserName"
$formSubmit.Range("D1").Value = "City"
$formSubmit.Range("E1").Value = "State"
$formSubmit.Range("F1").Value = "Zip"

$formSubmit.Range("F2").Value = 698473

[void]$formSubmit.Range("G7").Select()

# --- FormSubmitByParameter ---
$formSubmitByParameter.Activate()
$formSubmitByParameter.Range("A1").Value = "Expected String"
$formSubmitByParameter.Range("A2").Value = "Form has been submitted successfully!"
[void]$formSubmitByParameter.Range("A2").Select()

# --- FormSubmitByDataProvider (same two values as above) ---
$formSubmitByDataProvider.Activate()
$formSubmitByDataProvider.Range("A1").Value = "Expected String"
$formSubmitByDataProvider.Range("A2").Value = "Form has been submitted successfully!"
[void]$formSubmitByDataProvider.Range("A2").Select()

# Leave "FormSubmitByDataProvider" as the active tab.
$formSubmitByDataProvider.Activate()
